# Updated symbol list on Wed Dec 14 08:52:14 UTC 2022 with GitHub Actions
# Refreshes the Price (column D) for most rows and tweaks a couple of the
# Volume(1h) (column E) labels. Price values are stored as text in this
# sheet, so numeric-looking values are entered with a leading apostrophe
# to force text, then the cell style is reset to "Normal" so no stray
# quote-prefix formatting is left behind on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'276.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'22.98"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Value = "'6.382"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Value = "'0.06256"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Value = "'6.683"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9").Value = "'0.8309"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Value = "'0.1639"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Value = "'0.08371"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Value = "'0.03434"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Value = "'0.03103"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Value = "'0.09305"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Value = "'3.873"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Value = "'0.001636"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Value = "'0.04773"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "'0.006357"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'0.005558"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "19HotbitTokenHTB"
$ws.Range("D21").Value = "'0.001089"
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").Value = "'3.711"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Value = "'2.322"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'0.3375"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Value = "'0.1262"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Value = "'0.0002678"
$ws.Range("D27").Style = "Normal"
$ws.Range("D40").Value = "'0.04700"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Value = "'0.007023"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Value = "'0.1165"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.003348"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42CEJICEJIWorstin24h"
$ws.Range("D44").Value = "'0.01218"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006277"
$ws.Range("D45").Style = "Normal"
$ws.Range("D47").Value = "'0.8995"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.03495"
$ws.Range("D48").Style = "Normal"
